# "Back | update tests"
# 1) Compartments!E.. : a batch of rows get their "m_f_s_x type" dropdown
#    value flipped from "max" to "current" (shared-string index 192 -> 191).
# 2) GeneralCargo!A2:K2 : the placeholder "-" values are cleared out.
# 3) View state: selection/scroll position updated on GeneralCargo and
#    Compartments, and the workbook ends up with Compartments as the
#    active (visible) tab instead of Parameters.

$wb = $excel.ActiveWorkbook

# --- GeneralCargo: clear the placeholder row, move selection to M19 ----
# (Do this before touching Compartments so the *last* sheet activated by a
# Select() call is Compartments, matching the saved workbook's activeTab.)
$wsGeneralCargo = $wb.Worksheets.Item("GeneralCargo")
$wsGeneralCargo.Range("A2:K2").ClearContents()
$wsGeneralCargo.Range("M19").Select()

# --- Compartments: "max" -> "current" for the rows touched by this edit --
$wsCompartments = $wb.Worksheets.Item("Compartments")
$rowsToFlip = @(24,25,26,27,32,33,34,35,36,37,39,42,43,44,48,49,50,51,52)
foreach ($r in $rowsToFlip) {
    $wsCompartments.Range("E" + $r).Value = "current"
}

# Scroll/selection state + make Compartments the active sheet/tab.
$wsCompartments.Activate()
$wsCompartments.Range("G49:G50").Select()
